$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure contact-number-like text cells stay as text (preserve leading zeros)
$ws.Range("E2:E4").NumberFormat = "@"

$ws.Range("C2").Value = "AUTODSR_33662"
$ws.Range("D2").Value = "Father_D0A9"
$ws.Range("E2").Value = "03851754800"
$ws.Range("G2").Value = "EMPDD7C41"

$ws.Range("C3").Value = "AUTODSR_FC2B2"
$ws.Range("D3").Value = "Father_6326"
$ws.Range("E3").Value = "03854217700"
$ws.Range("G3").Value = "EMP7A5A9E"

$ws.Range("C4").Value = "AUTODSR_B0763"
$ws.Range("D4").Value = "Father_B318"
$ws.Range("E4").Value = "03855535900"
$ws.Range("G4").Value = "EMP706899"
